$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = "WbfuR431"
$ws.Range("B2").Value = 231106302
$ws.Range("C2").Value = "nhauqon87"
$ws.Range("D2").Value = "e38#VvE$"
$ws.Range("F2").Value = "ldhwuJFq"
$ws.Range("G2").Value = "fnJC"

# Row 3 updates
$ws.Range("A3").Value = "Kjyvo596"
$ws.Range("B3").Value = 231106301
$ws.Range("C3").Value = "njfvoeq85"
$ws.Range("D3").Value = "C9nh`$2%A"
$ws.Range("F3").Value = "aHYKTxAE"
$ws.Range("G3").Value = "mruj"
